# Update generated output values (想去人数 column F) on both the
# "展览" sheet and the aggregated "全部类型" sheet, matching the
# upstream data refresh captured in the commit diff.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 372
$wsExhibit.Range("F6").Value = 2000

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 372
$wsAll.Range("F10").Value = 2000
